$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Active" / "Yes" column in column E, matching the existing header/value layout
$ws.Range("E1").Value = "Active"
$ws.Range("E2").Value = "Yes"

# Update the selection to mirror the post-edit state (next empty row in the new column)
$ws.Range("E3").Select()
